# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

# Update the daily conversion summary text on Hoja1!A1
$nuevoTexto = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 13.79 = 55724.28 pesos`n✅ 55724.28 pesos = 13.73 = 962.84 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $nuevoTexto

# Update the rate figures on the "tasas" sheet
$wsTasas.Range("N10").Value = 72.5
$wsTasas.Range("O10").Value = 4040.01
$wsTasas.Range("N12").Value = 4060
$wsTasas.Range("O12").Value = 70.151
